# Writing of retest results in Excel file done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "Satut Bug" -> "Statut Bug"
$ws.Range("B3").Value = "Statut Bug"

# Helper: write a value as literal text (shared string), avoiding Excel's
# auto-conversion of number-/date-looking strings into numeric/date cells
# (and without leaving a stray NumberFormat-driven style behind).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Formula = "=""" + $val + """"
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# Retest results (rows 14-16)
Set-TextValue "A14" "1024"
Set-TextValue "B14" "KO"
Set-TextValue "C14" "2021-07-12"
Set-TextValue "D14" "16:54:45"

Set-TextValue "A15" "997"
Set-TextValue "B15" "KO"
Set-TextValue "C15" "2021-07-12"
Set-TextValue "D15" "16:55:27"

Set-TextValue "A16" "1011"
Set-TextValue "B16" "KO"
Set-TextValue "C16" "2021-07-12"
Set-TextValue "D16" "16:56:19"

$excel.CutCopyMode = $false
$ws.Range("B14").Select()
